$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("form")

# Update schedule data values
$ws.Range("H2").Value = "Преподаватель 1"
$ws.Range("F3").Value = "1"
$ws.Range("B5").Value = "2"
$ws.Range("B6").Value = "2"
$ws.Range("B7").Value = "2"

# Update the active cell selection shown in the sheet view
$ws.Activate()
$ws.Range("F6").Select()
